$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.4139
$ws.Range("B3").Value = 5.668499999999989
$ws.Range("D4").Value = -6.577799999999996
$ws.Range("B5").Value = 4.720400000000003
$ws.Range("D6").Value = -8.721299999999992
$ws.Range("C7").Value = -12.74110000000001
$ws.Range("A9").Value = -20.53629999999998
$ws.Range("C9").Value = -12.5162
$ws.Range("D10").Value = -6.389699999999996
$ws.Range("B11").Value = 5.387799999999996
$ws.Range("D11").Value = -8.314600000000002
$ws.Range("B12").Value = 5.130500000000001
$ws.Range("E12").Value = 11.8583
$ws.Range("A13").Value = -21.79780000000002
$ws.Range("A16").Value = -19.90519999999999
$ws.Range("E17").Value = 13.3956
$ws.Range("A18").Value = -22.47210000000001
$ws.Range("E19").Value = 13.3755
$ws.Range("A20").Value = -21.94880000000002
$ws.Range("B21").Value = 5.579399999999993
$ws.Range("C21").Value = -11.54579999999999
$ws.Range("D21").Value = -7.036999999999997
$ws.Range("E24").Value = 13.36669999999998
$ws.Range("D25").Value = -8.368799999999997
